# Update odds data for rows 2 and 3 (columns F:AO) as per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @{
    "F2" = 1.01;  "G2" = 1.01;  "H2" = 46;    "I2" = 1000;  "J2" = 150;
    "K2" = 990;   "L2" = 0;     "M2" = 0;     "N2" = 0;     "O2" = 0;
    "P2" = 1.98;  "Q2" = 2;     "R2" = 1.16;  "S2" = 6.6;   "T2" = 4.8;
    "U2" = 1.11;  "V2" = 1.01;  "W2" = 1.01;  "X2" = 1000;  "Y2" = 1000;
    "Z2" = 1000;  "AA2" = 1000; "AB2" = 1000; "AC2" = 1000; "AD2" = 1000;
    "AE2" = 1000; "AF2" = 2.06; "AG2" = 17;   "AH2" = 1000; "AI2" = 1000;
    "AJ2" = 5;    "AK2" = 55;   "AL2" = 1000; "AM2" = 1000; "AN2" = 32;
    "AO2" = 1000
}

$row3 = @{
    "F3" = 1.39;  "G3" = 1.41;  "H3" = 11.5;  "I3" = 13;    "J3" = 4.7;
    "K3" = 5;     "L3" = 0;     "M3" = 0;     "N3" = 5;     "O3" = 1.24;
    "P3" = 1.91;  "Q3" = 2.06;  "R3" = 1.29;  "S3" = 4.3;   "T3" = 1.92;
    "U3" = 1.96;  "V3" = 1.08;  "W3" = 3.4;   "X3" = 1000;  "Y3" = 1000;
    "Z3" = 1000;  "AA3" = 1000; "AB3" = 5.3;  "AC3" = 7.6;  "AD3" = 19.5;
    "AE3" = 980;  "AF3" = 5.8;  "AG3" = 8.4;  "AH3" = 23;   "AI3" = 790;
    "AJ3" = 13.5; "AK3" = 19;   "AL3" = 60;   "AM3" = 1000; "AN3" = 19.5;
    "AO3" = 980
}

foreach ($addr in $row2.Keys) {
    $ws.Range($addr).Value = $row2[$addr]
}

foreach ($addr in $row3.Keys) {
    $ws.Range($addr).Value = $row3[$addr]
}
